$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Cells.Item(12, 1).Value = 326292
$ws.Cells.Item(12, 2).Value = 'Fall Guys: Ultimate Knockout'
$ws.Cells.Item(12, 4).Value = '''2020-08-04'
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 'Fall Guys: Ultimate Knockout is a massively multiplayer online game developed by Mediatonic and published by Devolver Digital. Released on August 4, 2020, the game quickly gained popularity for its chaotic and comedic take on the battle royale genre. The game is available on multiple platforms, including Xbox One, PlayStation 4, PC, Nintendo Switch, and Xbox Series S/X, and supports cross-platform multiplayer, allowing players from different systems to compete against each other.
In Fall Guys: Ultimate Knockout, players are thrown into a series of increasingly challenging obstacle courses alongside up to 100 competitors online. The goal is to survive each round''s bizarre and whimsical challenges to ultimately emerge as the last contestant standing. The gameplay is characterized by its comically exaggerated physics, where characters bounce, tumble, and stumble their way through each course, adding a layer of unpredictability and humor to the competition. Players can customize their avatars with an array of quirky outfits, ranging from pineapple costumes to bunny hats, further enhancing the game''s lighthearted and playful atmosphere.
Upon release, Fall Guys was met with generally positive reviews, praised for its innovative approach to multiplayer gaming and its entertaining, family-friendly appeal. Critics highlighted the game''s colorful and cartoonish art style, which, as shown in the game''s artwork, adds to the overall fun and whimsical nature of the title. The game''s popularity was further fueled by its social media presence and streaming on platforms like Twitch, where it became a favorite among content creators and viewers alike. The game''s success led to it becoming one of the most downloaded titles on PlayStation Plus in its release month.
Culturally, Fall Guys: Ultimate Knockout made a significant impact by introducing a more accessible and less violent alternative to the traditional battle royale format. It resonated with a broad audience, including both casual and hardcore gamers, and became a social gaming phenomenon. Its influence extended beyond the gaming community, with various brands and celebrities expressing interest in collaborations and cross-promotions, further cementing its place in contemporary gaming culture.'
$ws.Cells.Item(12, 7).Value = '[''IGN. (2020). Fall Guys: Ultimate Knockout Review. Retrieved from https://www.ign.com/articles/fall-guys-ultimate-knockout-review'', ''Polygon. (2020). Fall Guys is the perfect game for when real life is too real. Retrieved from https://www.polygon.com/fall-guys-review'', ''The Verge. (2020). Fall Guys becomes the most downloaded PS Plus game of all time. Retrieved from https://www.theverge.com/2020/8/10/fall-guys-most-downloaded-ps-plus-game'', ''GameSpot. (2020). Fall Guys: Ultimate Knockout Review. Retrieved from https://www.gamespot.com/reviews/fall-guys-ultimate-knockout-review/1900-6417520/'', ''Mediatonic. (2020). Fall Guys: Ultimate Knockout Official Website. Retrieved from https://fallguys.com/'']'
$ws.Cells.Item(12, 8).Value = 'Average Rating: 3.74
Average Playtime: 6 hours
ESRB Rating: Everyone
Metacritic Score: 80
Platforms: Xbox One, PlayStation 4, PC, Nintendo Switch, Xbox Series S/X'
$ws.Cells.Item(12, 9).Value = '2025-03-27 19:14:46'
$ws.Cells.Item(12, 10).Value = 'https://media.rawg.io/media/games/5eb/5eb49eb2fa0738fdb5bacea557b1bc57.jpg'

# Row 13
$ws.Cells.Item(13, 1).Value = 274755
$ws.Cells.Item(13, 2).Value = 'Hades'
$ws.Cells.Item(13, 3).Value = 'Greg Kasavin, Amir Rao, Darren Korb, Jen Zee, Gavin Simon'
$ws.Cells.Item(13, 4).Value = '''2020-09-17'
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 'Hades is a rogue-like dungeon crawler developed and published by Supergiant Games, released on September 17, 2020. It is available on multiple platforms including PlayStation 5, Xbox Series S/X, PlayStation 4, Nintendo Switch, PC, and Xbox One. The game is notable for blending fast-paced action with deep storytelling and intricate character development, drawing upon the successful elements of Supergiant''s previous titles such as Bastion, Transistor, and Pyre.
In Hades, players assume the role of Zagreus, the immortal Prince of the Underworld, who is determined to escape the clutches of his father, Hades, the god of the dead. As Zagreus, players wield mythic weapons and receive powerful Boons from Olympian gods like Zeus, Athena, and Poseidon to enhance their abilities. The game is designed for high replayability, with each escape attempt offering a unique experience due to the ever-shifting Underworld and the various character builds players can discover. As shown in the game''s artwork, the rich, hand-painted environments and an original score by Darren Korb bring the atmospheric Underworld to life.
The development of Hades was led by Greg Kasavin, Amir Rao, Darren Korb, Jen Zee, and Gavin Simon. The game first entered early access in December 2018, allowing players to provide feedback that shaped its evolution. This iterative development process contributed significantly to its polished final release. The game features a fully-voiced cast and engaging narrative, with players able to befriend gods, ghosts, and monsters, thereby unlocking hundreds of unique story events.
Hades received widespread critical acclaim for its engaging gameplay, rich narrative, and audio-visual presentation. Reviews praised its deep replayability, character-driven storytelling, and the seamless integration of gameplay and narrative. It has been recognized with numerous awards, including the BAFTA Games Award for Best Game and Best Game Design in 2021. Its cultural impact is reflected not only in its critical success but also in its popularity among players, who have embraced its rich mythology and character depth.
The game''s success underscores the potential of indie games to deliver innovative and high-quality experiences, demonstrating how a small team can compete with larger studios by focusing on creativity and player engagement.'
$ws.Cells.Item(13, 7).Value = '[''IGN. (2020). Hades Review. Retrieved from https://www.ign.com/articles/hades-review'', ''Polygon. (2020). Hades is an Instant Classic. Retrieved from https://www.polygon.com/reviews/2020/9/17/hades-review'', "Kasavin, G. (2020). Supergiant''s Approach to Developing Hades. GDC Conference.", ''GameSpot. (2020). Why Hades is a Must-Play Game. Retrieved from https://www.gamespot.com/articles/why-hades-is-a-must-play-game/1100-6482730/'', ''Eurogamer. (2021). Hades Wins Big at the BAFTA Games Awards. Retrieved from https://www.eurogamer.net/hades-wins-big-at-bafta-games-awards'']'
$ws.Cells.Item(13, 8).Value = 'Average Rating: 4.44
Average Playtime: 10 hours
ESRB Rating: Teen
Metacritic Score: 93
Platforms: PlayStation 5, Xbox Series S/X, PlayStation 4, Nintendo Switch, PC, Xbox One'
$ws.Cells.Item(13, 9).Value = '2025-03-27 19:15:02'
$ws.Cells.Item(13, 10).Value = 'https://media.rawg.io/media/games/1f4/1f47a270b8f241e4676b14d39ec620f7.jpg'

# Row 14
$ws.Cells.Item(14, 1).Value = 22509
$ws.Cells.Item(14, 2).Value = 'Minecraft'
$ws.Cells.Item(14, 3).Value = 'Jens Bergensten, Daniel C418 Rosenfeld, Henrik "carnalizer" Pettersson, Markus Persson, Markus "Junkboy" Toivonen, Jon Kågström, Aron Nieminen, Mattis "Bomb Boy" Grahm, Daniel Frisk'
$ws.Cells.Item(14, 4).Value = '''2009-05-10'
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(14, 6).Value = 'Minecraft is an iconic sandbox video game that was initially released on May 10, 2009. Developed by a team that included notable figures such as Markus Persson and Jens Bergensten, and published by Mojang and Microsoft Studios, Minecraft has become one of the most influential and best-selling games of all time. Available on multiple platforms including Android, iOS, PlayStation, Xbox, and PC, the game is renowned for its open-world design and limitless creative potential.
The core gameplay of Minecraft revolves around mining resources and crafting items to build various structures. The game''s world is procedurally generated and composed of 3D cubes representing different materials and terrains, allowing players to interact with and modify their environment extensively. Players can choose from multiple game modes, including Survival, Creative, and Hardcore. Survival mode challenges players to gather resources, hunt for food, and defend against monsters, while Creative mode offers unlimited resources for unfettered construction and exploration. Hardcore mode introduces permanent death, adding a layer of challenge for seasoned players. Additionally, a multiplayer mode facilitates shared experiences and competitive player-versus-player interactions.
Minecraft is distinguished by its crude, pixelated visual style, reminiscent of Lego blocks, which has become a cultural icon in its own right. The game''s aesthetic simplicity belies its complex and engaging gameplay mechanics, which have inspired a wide array of myths and fan fiction, such as the legendary but non-existent character Herobrine.
The reception of Minecraft has been overwhelmingly positive, with a rating of 4.43. Critics and players alike have praised its innovative approach to gaming, its vast creative possibilities, and its engaging, community-driven multiplayer aspect. The game has also had a significant cultural impact, influencing numerous aspects of popular culture and being integrated into educational settings to teach concepts ranging from mathematics to environmental stewardship.
As shown in the game''s artwork, Minecraft''s visual style and gameplay mechanics have made it a staple in the gaming world, fostering a dedicated community of modders and players who continue to expand and shape the game. The official website, classic.minecraft.net, offers an entry point into the expansive world of Minecraft, a game that continues to evolve and captivate new generations of players.'
$ws.Cells.Item(14, 7).Value = '[''Persson, M. (2009). *Minecraft*. Mojang.'', ''Totilo, S. (2011). *The Minecraft phenomenon*. Kotaku. Retrieved from https://kotaku.com/the-minecraft-phenomenon-5820413'', ''IGN Staff. (2011). *Minecraft Review*. IGN. Retrieved from https://www.ign.com/articles/2011/11/24/minecraft-review'', ''Schreier, J. (2015). *How Minecraft became the best-selling PC game of all time*. Kotaku. Retrieved from https://kotaku.com/how-minecraft-became-the-best-selling-pc-game-of-all-ti-1696666722'', ''GameSpot Staff. (2010). *Minecraft: Game Overview*. GameSpot. Retrieved from https://www.gamespot.com/games/minecraft/'']'
$ws.Cells.Item(14, 8).Value = 'Average Rating: 4.43
Average Playtime: 26 hours
ESRB Rating: Everyone 10+
Metacritic Score: 83
Platforms: Android, PS Vita, PlayStation 4, PlayStation 3, Xbox 360, Nintendo 3DS, Nintendo Switch, macOS, PC, iOS, Wii U, Xbox One, Linux'
$ws.Cells.Item(14, 9).Value = '2025-03-27 19:15:51'
$ws.Cells.Item(14, 10).Value = 'https://media.rawg.io/media/games/b4e/b4e4c73d5aa4ec66bbf75375c4847a2b.jpg'
